# Weekly update: insert the new week's price row for
# "Hortaliza, Mercado Mayorista Lo Valledor de Santiago - Ciboulette"
# at row 260, pushing the existing rows (260..380) down by one
# (new last row becomes 381).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 260; everything below shifts down.
$ws.Rows.Item(260).Insert()

# Populate the newly inserted row with this week's data point.
$ws.Cells.Item(260, 1).Value  = 6
$ws.Cells.Item(260, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(260, 3).Value  = "Metropolitana"
$ws.Cells.Item(260, 4).Value  = 44609
$ws.Cells.Item(260, 5).Value  = 13
$ws.Cells.Item(260, 6).Value  = 100112039
$ws.Cells.Item(260, 7).Value  = "Ciboulette"
$ws.Cells.Item(260, 8).Value  = "Sin especificar"
$ws.Cells.Item(260, 9).Value  = "Primera"
$ws.Cells.Item(260, 10).Value = 930
$ws.Cells.Item(260, 11).Value = 900
$ws.Cells.Item(260, 12).Value = 1000
$ws.Cells.Item(260, 13).Value = 947
$ws.Cells.Item(260, 14).Value = "`$/docena de atados"
$ws.Cells.Item(260, 15).Value = "Región Metropolitana"
$ws.Cells.Item(260, 16).Value = 316
$ws.Cells.Item(260, 17).Value = 3
$ws.Cells.Item(260, 18).Value = "Hortaliza"

# Keep the date cell formatted like the rest of column D.
$ws.Cells.Item(260, 4).NumberFormat = $ws.Cells.Item(261, 4).NumberFormat
